$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "26/5/2015"
$ws.Range("B10").Value = "Comment for the demo project: Get list comment of a book and add new comment on client"

$ws.Range("B10").Select()
